$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A37").Value = "USC2X3"
$ws.Range("B37").Value = "Cuchillo de limpieza Ricoh"
$ws.Range("C37").Value = "Aficio 1013, 120, 1515, 171, MP161 MP201 MP301"
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 100000
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 0
$ws.Range("H37").Formula = "=(E37-D37)*G37"
$ws.Range("I37").Formula = "=D37*F37"
$ws.Range("J37").Value = 0

$wb.Save()
